$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Sheet1 (quality_comparison) ---

# C1: reset to default, then apply top+bottom border (-> new style reusing borderId 4)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1: start from C1's finished top+bottom style, then add the right edge
# (-> new style reusing borderId 5), avoiding revisiting the "top only" state
$d1 = $ws1.Range("D1")
$c1.Copy()
$d1.PasteSpecial(-4122)  # xlPasteFormats
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight

# reuse the exact same computed styles for sheet2's matching cells via copy/paste of formats,
# instead of rebuilding border-by-border (avoids stray intermediate styles being recorded)
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# C2: "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet2 (computational_comparison) ---
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 previously held an empty placeholder cell; remove it entirely
$ws2.Range("G5").ClearContents()

Write-Host "done"
